$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in additional hours for Tuesday (C) and Wednesday (D) columns for rows 2-7
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5

$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 6

$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 6

$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 6

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1

# Update view selection
$ws.Range("F10").Select()
